$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.281.35"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.440.26"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.77"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.44"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").Value = "3.431.79"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -2.21%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  +0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.07"
$ws.Range("E11").Value = "  -3.13%  "
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.20"
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000269"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "3.994.75"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "581.20"
$ws.Range("E17").Value = "  -4.76%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.445.72"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "69.352.32"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.12"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.85"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "95.73"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.14"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.69"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.60"
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("E31").Value = "  -3.31%  "
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  -6.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.53"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "574.94"
$ws.Range("E35").Value = "  -9.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.50"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0474"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "55.72"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.14"
$ws.Range("E42").Value = "  -11.60%  "
$ws.Range("D43").Value = "3.238.02"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("D44").Value = "0.0₃0682"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.11"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("E48").Value = "  -5.74%  "
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.88"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  -0.04%  "
